# Apply "back-of-a-napkin calc" relay-settings update to the F3 relay sheet.
#
# For rows 2-10 (the 9 relay rows):
#   - Column G (pickup current, "Ipickup") is bumped to a new, coarser value
#     and picks up an explicit centered-alignment style (no longer just the
#     generic numeric style).
#   - Column K (time dial) gets a computed value (2 decimal place format) and
#     a new style: centered alignment + "0.00" number format.
#   - Column N / O (min/max multiplier band) are tightened from 0.9/1.1 to
#     0.8/1.2, values only - the style stays the same.
# The active selection also moves from E11 to D15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pickup current (G) and computed time-dial (K) values per relay row.
$G = @{
    2 = 1200
    3 = 1200
    4 = 1200
    5 = 1200
    6 = 1200
    7 = 1200
    8 = 3500
    9 = 3500
    10 = 3500
}

$K = @{
    2 = 0.80187537387448016
    3 = 0.20046884346862004
    4 = 0.90210979560879001
    5 = 0.50117210867155004
    6 = 0.40093768693724008
    7 = 0.40093768693724008
    8 = 1.1856300170858387
    9 = 1.1856300170858387
    10 = 2.4700625355954968
}

foreach ($row in 2..10) {
    # --- Column G: new pickup current, centered alignment ---
    $gCell = $ws.Cells.Item($row, 7)
    $gCell.Value = $G[$row]
    $gCell.HorizontalAlignment = -4108   # xlCenter

    # --- Column K: computed time dial, centered + 2-decimal format ---
    $kCell = $ws.Cells.Item($row, 11)
    $kCell.NumberFormat = "0.00"
    $kCell.HorizontalAlignment = -4108   # xlCenter
    $kCell.Value = $K[$row]

    # --- Column N / O: tighten the multiplier band (values only) ---
    $ws.Cells.Item($row, 14).Value = 0.8
    $ws.Cells.Item($row, 15).Value = 1.2
}

# Move the active selection from E11 to D15.
$ws.Range("D15").Select()
